# Update the "想去人数" (interested-count) figures in column F across the
# workbook's sheets to reflect the latest generated data snapshot.
# Sheet order: 1 = 展览, 2 = 演出, 3 = 本地生活, 4 = 全部类型

$wb = $excel.ActiveWorkbook

# Sheet 1 - 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 555
$ws1.Range("F6").Value  = 1613
$ws1.Range("F9").Value  = 748
$ws1.Range("F10").Value = 2717
$ws1.Range("F12").Value = 1823
$ws1.Range("F14").Value = 304
$ws1.Range("F15").Value = 710
$ws1.Range("F16").Value = 5
$ws1.Range("F17").Value = 6246
$ws1.Range("F18").Value = 237
$ws1.Range("F19").Value = 89
$ws1.Range("F21").Value = 3398
$ws1.Range("F22").Value = 881
$ws1.Range("F26").Value = 2467
$ws1.Range("F28").Value = 379
$ws1.Range("F32").Value = 1316
$ws1.Range("F35").Value = 87
$ws1.Range("F38").Value = 1497
$ws1.Range("F39").Value = 30
$ws1.Range("F40").Value = 1457

# Sheet 2 - 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F16").Value = 147
$ws2.Range("F18").Value = 268

# Sheet 3 - 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 924
$ws3.Range("F4").Value = 259
$ws3.Range("F6").Value = 45

# Sheet 4 - 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 924
$ws4.Range("F5").Value  = 259
$ws4.Range("F6").Value  = 555
$ws4.Range("F7").Value  = 45
$ws4.Range("F16").Value = 2717
$ws4.Range("F21").Value = 1823
$ws4.Range("F24").Value = 304
$ws4.Range("F25").Value = 710
$ws4.Range("F26").Value = 6246
$ws4.Range("F27").Value = 237
$ws4.Range("F28").Value = 89
$ws4.Range("F30").Value = 3398
$ws4.Range("F31").Value = 881
$ws4.Range("F36").Value = 2467
$ws4.Range("F37").Value = 379
$ws4.Range("F39").Value = 1316
$ws4.Range("F41").Value = 268
$ws4.Range("F45").Value = 87
$ws4.Range("F48").Value = 30
$ws4.Range("F50").Value = 1457
